$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 318
$ws1.Range("F4").Value = 1278
$ws1.Range("F5").Value = 631

# Sheet "全部类型" (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 318
$ws4.Range("F4").Value = 1278
$ws4.Range("F6").Value = 631
